# The commit adds one new price-record row for "Perejil" at the top of the
# weekly data block (row 76), pushing all the existing records (old rows
# 76-164) down by one row (new rows 77-165). This mirrors the "Fruta /
# hortaliza, semanal" weekly-refresh pattern used across these workbooks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 76; Excel shifts rows 76:164 down to 77:165
# and the sheet dimension grows from R164 to R165 automatically.
$ws.Rows("76:76").Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(76, 1).Value = 8
$ws.Cells.Item(76, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(76, 3).Value = "Coquimbo"
$ws.Cells.Item(76, 4).Value = 44803
$ws.Cells.Item(76, 5).Value = 4
$ws.Cells.Item(76, 6).Value = 100112044
$ws.Cells.Item(76, 7).Value = "Perejil"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 2400
$ws.Cells.Item(76, 11).Value = 2000
$ws.Cells.Item(76, 12).Value = 2500
$ws.Cells.Item(76, 13).Value = 2250
$ws.Cells.Item(76, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(76, 16).Value = 1500
$ws.Cells.Item(76, 17).Value = 1.5
$ws.Cells.Item(76, 18).Value = "Hortaliza"
